$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.129.66'
$ws.Range("E2").Value = '  +3.41%  '
$ws.Range("D3").Value = '2.982.97'
$ws.Range("E3").Value = '  +2.32%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '597.06'
$ws.Range("E5").Value = '  +1.14%  '
$ws.Range("E6").Value = '  +1.10%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '2.981.63'
$ws.Range("E8").Value = '  +2.39%  '
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("D10").Value = '7.43'
$ws.Range("E10").Value = '  +7.03%  '
$ws.Range("E11").Value = '  +2.79%  '
$ws.Range("E12").Value = '  +3.30%  '
$ws.Range("E13").Value = '  +5.57%  '
$ws.Range("D14").Value = '33.76'
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("D16").Value = '3.476.04'
$ws.Range("E16").Value = '  +2.27%  '
$ws.Range("D17").Value = '62.948.04'
$ws.Range("E17").Value = '  +3.27%  '
$ws.Range("E18").Value = '  +1.49%  '
$ws.Range("D19").Value = '2.984.85'
$ws.Range("E19").Value = '  +2.42%  '
$ws.Range("D20").Value = '444.91'
$ws.Range("E20").Value = '  +2.62%  '
$ws.Range("D21").Value = '13.63'
$ws.Range("E21").Value = '  +2.06%  '
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("D23").Value = '7.17'
$ws.Range("E23").Value = '  +1.13%  '
$ws.Range("D24").Value = '82.37'
$ws.Range("E24").Value = '  +1.19%  '
$ws.Range("D25").Value = '10.99'
$ws.Range("E25").Value = '  +1.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.20'
$ws.Range("E26").Value = '  +3.61%  '
$ws.Range("D27").Value = '2.17'
$ws.Range("E27").Value = '  -1.57%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("E29").Value = '  +1.52%  '
$ws.Range("E31").Value = '  -5.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.70'
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = '0.0₃0890'
$ws.Range("E35").Value = '  +2.57%  '
$ws.Range("D36").Value = '0.997'
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.70'
$ws.Range("E37").Value = '  +1.59%  '
$ws.Range("D38").Value = '2.07'
$ws.Range("E38").Value = '  +4.57%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '49.84'
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").Value = '2.99'
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("E41").Value = '  +1.24%  '
$ws.Range("E42").Value = '  -2.38%  '
$ws.Range("D43").Value = '0.287'
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("D44").Value = '39.04'
$ws.Range("E44").Value = '  -6.63%  '
$ws.Range("D45").Value = '374.52'
$ws.Range("E45").Value = '  -0.42%  '
$ws.Range("D46").Value = '0.0346'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("D47").Value = '2.710.92'
$ws.Range("E47").Value = '  +0.81%  '
$ws.Range("D48").Value = '135.04'
$ws.Range("E48").Value = '  +1.61%  '
$ws.Range("D50").Value = '23.48'
$ws.Range("E50").Value = '  -1.36%  '
$ws.Range("E51").Value = '  +0.14%  '
